$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-21 Wednesday" "2024-02-22 Thursday"

Replace-Text "48÷5=" "93÷2="
Replace-Text "37÷8=" "67÷8="
Replace-Text "10÷4=" "60÷9="
Replace-Text "62÷2=" "75÷7="
Replace-Text "30÷9=" "53÷8="
Replace-Text "10÷3=" "59÷2="
Replace-Text "98÷7=" "39÷4="
Replace-Text "92÷8=" "45÷8="
Replace-Text "72÷6=" "87÷8="
Replace-Text "36÷9=" "71÷9="
Replace-Text "10÷7=" "23÷5="
Replace-Text "85÷3=" "82÷8="
Replace-Text "10÷2=" "24÷4="
Replace-Text "97÷9=" "89÷2="
Replace-Text "51÷7=" "79÷3="
Replace-Text "35÷2=" "97÷3="
Replace-Text "36÷4=" "90÷8="
Replace-Text "23÷4=" "60÷3="
Replace-Text "44÷2=" "91÷3="
Replace-Text "24÷3=" "80÷7="
Replace-Text "46÷6=" "48÷4="
Replace-Text "76÷5=" "62÷4="
Replace-Text "92÷5=" "85÷6="
Replace-Text "33÷9=" "65÷9="
Replace-Text "26÷4=" "21÷6="
